$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52; this shifts existing rows 52-102 down to 53-103,
# preserving all of their data/formatting (matches the target diff exactly).
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record's data.
$ws.Cells.Item(52, 1).Value = 1
$ws.Cells.Item(52, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value = 45174
$ws.Cells.Item(52, 5).Value = 15
$ws.Cells.Item(52, 6).Value = 100112009
$ws.Cells.Item(52, 7).Value = "Acelga"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 180
$ws.Cells.Item(52, 11).Value = 1800
$ws.Cells.Item(52, 12).Value = 2000
$ws.Cells.Item(52, 13).Value = 1889
$ws.Cells.Item(52, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(52, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value = 630
$ws.Cells.Item(52, 17).Value = 3
$ws.Cells.Item(52, 18).Value = "Hortaliza"
